$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (existing row, was "some wires") becomes "switch" + quantity, with a
#     plain (non-shared) formula, matching how Excel breaks a shared-formula
#     group when the sheet gets restructured around it ---

# First, remove the old trailing blank formula row (row 18) that isn't needed
# in the final layout, then insert 5 fresh rows at row 15 to make room for the
# new BOM entries (this pushes "some wires" / "you may also place..." / the
# blank row down to rows 20-22).
$ws.Rows("18:18").Delete()
$ws.Rows("15:19").Insert()

# The row-insert copies the formatting of row 14 across the whole row width;
# drop the C/D formatting carried into the new rows since those columns stay
# empty for this block.
$ws.Range("C15:D19").Clear()

# Re-apply the plain body style (same as column A/B) to the E column for the
# freshly inserted rows.
$ws.Range("A14").Copy()
$ws.Range("E15:E19").PasteSpecial(-4122)

# New BOM rows. Shared-string entries are appended in first-write order, so
# the values are typed in the same order the author must have used (switch,
# then the M3 x 30 screws, then nuts/washers/M3 x 16 screws) to line up with
# the resulting shared string table.
$ws.Range("A15").Value = "switch"
$ws.Range("B15").Value = 1
$ws.Range("E15").Formula = "=B15*C15"

$ws.Range("A19").Value = "m3 x 30 screw"
$ws.Range("B19").Value = 8

$ws.Range("A16").Value = "m3 nuts"
$ws.Range("B16").Value = 12

$ws.Range("A17").Value = "m3 washers"
$ws.Range("B17").Value = 12

$ws.Range("A18").Value = "m3 x 16 screw"
$ws.Range("B18").Value = 4

# Rows 20-22 are the old "some wires" / "you may also place..." / blank rows,
# shifted down by the insert above. Clear their leftover shared-formula
# remnants in column E and rebuild the exact target cell layout.
$ws.Range("E20:E22").ClearContents()
$ws.Range("E20").Clear()

$ws.Range("A14").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("F20").PasteSpecial(-4122)

# Move the active selection to E21, matching the author's final cursor spot.
$ws.Range("E21").Select()
